# Swap the presentation's theme palette from "Integral" to the stock
# "Office Theme" colour scheme (the slide master's theme -- ppt/theme/theme1.xml
# -- and, through it, the notes/handout master which shares the same
# in-memory Theme object).
#
# PowerPoint COM represents each of the 12 theme colour slots via
# Theme.ThemeColorScheme.Item(n).RGB, in the fixed order:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
#   8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink
# RGB values use the standard OLE COLORREF packing: R + (G*256) + (B*65536).

$p = $ppt.ActivePresentation
$sm = $p.SlideMaster
$theme = $sm.Theme
$colors = $theme.ThemeColorScheme

function RGBVal($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

# Office Theme colour scheme (dk1..folHlink)
$officeColors = @(
    (RGBVal 0x00 0x00 0x00),   # 1  dk1      000000
    (RGBVal 0xFF 0xFF 0xFF),   # 2  lt1      FFFFFF
    (RGBVal 0x44 0x54 0x6A),   # 3  dk2      44546A
    (RGBVal 0xE7 0xE6 0xE6),   # 4  lt2      E7E6E6
    (RGBVal 0x5B 0x9B 0xD5),   # 5  accent1  5B9BD5
    (RGBVal 0xED 0x7D 0x31),   # 6  accent2  ED7D31
    (RGBVal 0xA5 0xA5 0xA5),   # 7  accent3  A5A5A5
    (RGBVal 0xFF 0xC0 0x00),   # 8  accent4  FFC000
    (RGBVal 0x44 0x72 0xC4),   # 9  accent5  4472C4
    (RGBVal 0x70 0xAD 0x47),   # 10 accent6  70AD47
    (RGBVal 0x05 0x63 0xC1),   # 11 hlink    0563C1
    (RGBVal 0x95 0x4F 0x72)    # 12 folHlink 954F72
)

for ($i = 1; $i -le $colors.Count; $i++) {
    $colors.Item($i).RGB = $officeColors[$i - 1]
}
